# Update the "name" column text for two rows in the "metric" sheet:
#  - sandwich-population row: "Sandwich Generation Population" -> "Sandwich Generation: Population"
#  - sandwich-time row:       "Daily Caregiving Hours"         -> "Sandwich Generation: Daily Caregiving Hours"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metric")

$ws.Range("C4").Value = "Sandwich Generation: Population"
$ws.Range("C5").Value = "Sandwich Generation: Daily Caregiving Hours"

# Widen column C so the longer labels still fit (matches the recorded column width change;
# 40.67 is the input that lands on the nearest achievable stored width to 41.42578125).
$ws.Columns.Item(3).ColumnWidth = 40.67

# Restore the active cell/selection on the metric sheet to C5 (matches the saved view state).
$ws.Activate()
$ws.Range("C5").Select()
